$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the data block (rows 3-15) directly with the final values from the
# updated asset list documentation, rather than shifting rows (the trailing
# rows 17-19 must stay put).

# Row 3: button-hover (unchanged)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "button-hover"
$ws.Range("C3").Value = "whoosh noise, similar frequency throughout sound but increases a bit, alerts user they are hovering an option"
$ws.Range("D3").Value = "Interface"
$ws.Range("F3").Value = "implemented in Unity"

# Row 4: button-click (unchanged)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "button-click"
$ws.Range("C4").Value = "glassy-sounding, high to low pitch curve, alerts user they clicked a button"
$ws.Range("D4").Value = "Interface"
$ws.Range("F4").Value = "implemented in Unity"

# Row 5: beach-waves (description trimmed)
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "beach-waves"
$ws.Range("C5").Value = "crashing waves on the beach, relaxing mood, will be used as background ambience "
$ws.Range("D5").Value = "Ambient "
$ws.Range("F5").Value = "implemented in Unity"

# Row 6: background-music (unchanged)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "background-music"
$ws.Range("C6").Value = "melodic seamless loop of bright-sounding music"
$ws.Range("D6").Value = "Music"
$ws.Range("F6").Value = "implemented in Unity"

# Row 7: gravel-walk-param (renamed, now crab x specific)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "gravel-walk-param(footsteps1,2,3)"
$ws.Range("C7").Value = "whenever player assets move, they will be making a seamless sand moving sound, sound of crab walking on beach(not realistic), meant for crab x"
$ws.Range("D7").Value = "Sound Effect"
$ws.Range("F7").Value = "implemented in Unity"

# Row 8: gravel-walk-param2 (new, crab y specific)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "gravel-walk-param2(footsteps1,2,3)"
$ws.Range("C8").Value = "whenever player assets move, they will be making a seamless sand moving sound, sound of crab walking on beach(not realistic), meant for crab y"
$ws.Range("D8").Value = "Sound Effect"
$ws.Range("F8").Value = "implemented in Unity"

# Row 9: move-object (shifted down from old row 8)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "move-object"
$ws.Range("C9").Value = "planned to sound like heavy object moving across sand/gravel, alerts user object is currently being pushed"
$ws.Range("D9").Value = "Sound Effect"
$ws.Range("F9").Value = "implemented in Unity"

# Row 10: button-activated (shifted down, status now incomplete)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "button-activated"
$ws.Range("C10").Value = "a very clear click sound will be played when button is pressed"
$ws.Range("D10").Value = "Sound Effect"
$ws.Range("F10").Value = "in Unity but incomplete"

# Row 11: gate-open (shifted down)
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "gate-open"
$ws.Range("C11").Value = "sounds like a metal latch unlocking or moving, will be played when player opens the gate on top of the exit"
$ws.Range("D11").Value = "Sound Effect"
$ws.Range("F11").Value = "in Unity but incomplete"

# Row 12: gate-close (shifted down)
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "gate-close"
$ws.Range("C12").Value = "sounds like a metal latch unlocking or moving but will sound pitched down, will be played gate is opened, but then closes"
$ws.Range("D12").Value = "Sound Effect"
$ws.Range("F12").Value = "implemented in Unity"

# Row 13: level-complete (shifted down)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "level-complete"
$ws.Range("C13").Value = "bright sounding artificial excitement noise, sound is played when crab enters ungated hole and progresses to next level"
$ws.Range("D13").Value = "Sound Effect"
$ws.Range("F13").Value = "implemented in Unity"

# Row 14: button-unactivated (shifted down, new row, no status)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "button-unactivated"
$ws.Range("C14").Value = "pitched down click sound will be played when button is unpressed"
$ws.Range("D14").Value = "Sound Effect"

# Row 15: wall collision (brand-new row, no status)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "wall collision"
$ws.Range("C15").Value = "small bump noise"
$ws.Range("D15").Value = "Sound Effect"

# View/selection tidy-up: drop the frozen topLeftCell and move the active selection to C8
$ws.Range("C8").Select()
